$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header row height changes from 40.5 to 33.75 ---
$ws.Rows.Item(1).RowHeight = 33.75

# --- New column J (year 2021) values, rows 4-23 ---
# Each cell's number format/style mirrors the corresponding cell in column H
# (the 2019 column), which is how the original author extended the table.

$yearValues = @{
    4  = 2021
    5  = 1006091.2
    6  = 2092.6999999999998
    7  = 211904.6
    8  = 228945.8
    9  = 6780.6
    10 = 92.5
    11 = 9456.7999999999993
    12 = 92470.9
    13 = 656.4
    14 = 3692
    15 = 59559.1
    16 = 53592.2
    17 = 11799.2
    18 = 316755
    19 = 901
    20 = 76.5
    21 = 1672.3
    22 = 5539.9
    23 = 103.7
}

foreach ($row in 4..23) {
    # Copy formatting from column H (same row) into column J, then set the value.
    $ws.Cells.Item($row, 8).Copy()
    $ws.Cells.Item($row, 10).PasteSpecial(-4122)
    $ws.Cells.Item($row, 10).Value2 = $yearValues[$row]
}

# Row 22's H column cell never had its number format actually wired up
# (applyNumberFormat was set but numFmtId stayed General), so the new J22
# cell needs the real "#,##0.0" numeric format applied on top of the copied
# look of I22 (which already carries applyNumberFormat="1").
$ws.Cells.Item(22, 9).Copy()
$ws.Cells.Item(22, 10).PasteSpecial(-4122)
$ws.Cells.Item(22, 10).NumberFormat = "#,##0.0"
$ws.Cells.Item(22, 10).Value2 = 5539.9

$excel.CutCopyMode = 0

# --- Update the active selection shown when the sheet was last saved ---
$ws.Range("K3").Select()
